$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 107
$ws.Range("B7").Value = "Shakib"
$ws.Range("C7").Value = "Shak Forid"
$ws.Range("D7").Value = "Aklima"
$ws.Range("E7").Value = 677287
$ws.Range("F7").Value = "Male"
$ws.Range("G7").Value = "Computer Science & Technology (85)"
$ws.Range("H7").Value = "23/06/2003"
$ws.Range("I7").Value = "Habiganj Polytechnic Institute (63010)"
$ws.Range("J7").Value = "Gopaya"
$ws.Range("K7").Value = "Habiganj Sadar"
$ws.Range("L7").Value = "Habiganj"
$ws.Range("M7").Value = "2021-22"
